$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$esc = [char]0x1B

# Replace the "bat_monitor:" keyword row (row 12) with the new "cpu_start:" keyword/example
$ws.Range("C12").Value = "cpu_start:"
$ws.Range("D12").Value = "[11:39:21.922]IN¡û¡ô${esc}[0;32mI (92) cpu_start: Pro cpu start user code${esc}"

# Update the selection to match the saved file (single cell C12 selected)
$ws.Range("C12").Select() | Out-Null
